$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(29).Insert()

$ws.Cells.Item(29, 1).Value = 3
$ws.Cells.Item(29, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(29, 3).Value = "Coquimbo"
$ws.Cells.Item(29, 4).Value = 44560
$ws.Cells.Item(29, 5).Value = 5
$ws.Cells.Item(29, 6).Value = 100112030
$ws.Cells.Item(29, 7).Value = "Poroto granado"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 35
$ws.Cells.Item(29, 11).Value = 30000
$ws.Cells.Item(29, 12).Value = 30000
$ws.Cells.Item(29, 13).Value = 30000
$ws.Cells.Item(29, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(29, 16).Value = 1200
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
